$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mythril (row 3) and Semgrep (row 4) tool-capacity flags flipped from 0 to 1
# after the Mythril v0.24.8 upgrade.
$ws.Range("B3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 1

# Refresh the view: zoom to 100% and move the selection to C12.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("C12").Select()
